$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Replace the old report-generation date (2024-02-05) with the new one
# (2024-12-02) in every path stored in column C, for every data row.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 3).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    $value = $cell.Value2
    if ($value -ne $null -and $value -like "*2024-02-05*") {
        $cell.Value2 = $value -replace "2024-02-05", "2024-12-02"
    }
}
